$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'25.067.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").Value = "'1.651.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.37%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'237.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.72%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -6.49%  "
$ws.Range("D8").Value = "'0.2617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.19%  "
$ws.Range("D9").Value = "'0.06015"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.07187"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'1.651.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("D12").Value = "'14.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("D13").Value = "'0.6221"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("D14").Value = "'4.596"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "'73.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.88%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'0.9996"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'25.056.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.66%  "
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").Value = "'0.000006625"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("E21").Value = "  +5.48%  "
$ws.Range("D22").Value = "'1.863.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.21%  "
$ws.Range("D23").Value = "'8.628"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").Value = "'5.292"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Value = "'131.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").Value = "'15.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Value = "'1.403"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.13%  "
$ws.Range("D28").Value = "'103.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "'1.685"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.76%  "
$ws.Range("D30").Value = "'3.784"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("D31").Value = "'0.07911"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("D32").Value = "'3.610"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").Value = "'0.04610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").Value = "'2.593"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "'0.9427"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.84%  "
$ws.Range("D36").Value = "'0.5768"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.79%  "
$ws.Range("D37").Value = "'2.604"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.47%  "
$ws.Range("D38").Value = "'0.01561"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("D39").Value = "'0.9995"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.8173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.08%  "
$ws.Range("D41").Value = "'1.831"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.17%  "
$ws.Range("D42").Value = "'98.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").Value = "'0.3733"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").Value = "'4.798"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.56%  "
$ws.Range("D45").Value = "'0.1144"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "'6.114"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("D47").Value = "'0.05186"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "'29.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "'51.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.99%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'0.3337"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.00%  "
